$wb = $excel.ActiveWorkbook

# Capture header style source (existing header cell with bold+border+center formatting)
$styleSource = $wb.Worksheets.Item("Az_Asctb_cts_perfect_matches").Range("A1")

# Add the new worksheet as the last sheet in the workbook
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "Final_Matches"

# Gold standard mapping comparison data (AZ.CT/LABEL vs ASCTB.CT/LABEL)
$data = @(
    ,("AZ.CT/LABEL", "ASCTB.CT/LABEL")
    ,("club cell", "club cell of bronchiole")
    ,("myofibroblast cell", "secondary crest myofibroblasts")
    ,("B cell", "B cell")
    ,("natural killer cell", "natural killer")
    ,("pericyte cell", "pericyte")
    ,("plasmacytoid dendritic cell", "plasmacytoid dendritic cell")
    ,("plasma cell", "plasma cell")
    ,("mature NK T cell", "NK T cell")
    ,("classical monocyte", "Classical Monocyte")
    ,("non-classical monocyte", "Non-classical monocyte")
    ,("naive thymus-derived CD4-positive, alpha-beta T cell", "CD4+ T cell naive")
    ,("naive thymus-derived CD8-positive, alpha-beta T cell", "CD8+ T cell naive")
    ,("type I pneumocyte", "type I pneumocyte")
    ,("type II pneumocyte", "type II pneumocyte")
    ,("CD4-positive, alpha-beta T cell", "T cell")
    ,("CD8-positive, alpha-beta T cell", "T cell")
    ,("CD14-positive, CD16-negative classical monocyte", "classical monocyte")
    ,("CD14-low, CD16-positive monocyte", "non-classical monocyte")
    ,("ciliated cell", "ciliated columnar cell of tracheobronchial tree")
    ,("serous secreting cell", "serous cell of epithelium of trachea")
    ,("smooth muscle cell", "smooth muscle cell of trachea")
    ,("vascular associated smooth muscle cell", "blood vessel smooth muscle cell")
    ,("bronchial epithelial cell", "ciliated cell of the bronchus")
    ,("serous secreting cell", "serous cell of epithelium of bronchus")
    ,("smooth muscle cell", "bronchial smooth muscle cell")
    ,("vascular associated smooth muscle cell", "smooth muscle cell of pulmonary artery")
    ,("fibroblast", "fibroblast of pulmonary artery")
    ,("smooth muscle cell", "tracheobronchial smooth muscle cell")
    ,("serous secreting cell", "serous cell of epithelium of lobular bronchiole")
    ,("serous secreting cell", "serous cell of epithelium of terminal bronchiole")
    ,("vascular associated smooth muscle cell", "lymphatic vessel smooth muscle cell")
    ,("macrophage", "Alveolar Macrophage")
    ,("myeloid leukocyte", "Basophil")
    ,("myeloid leukocyte", "lung parenchyma resident eosinophil")
    ,("myeloid leukocyte", "neutrophil")
    ,("lymphocyte", "T cell")
    ,("plasmacytoid dendritic cell", "myeloid dendritic cell cDC1")
    ,("myeloid dendritic cell", "myeloid dendritic cell cDC2")
    ,("lymphocyte", "regulatory T cell")
    ,("lymphocyte", "CD4+ T cell central memory")
    ,("lymphocyte", "CD8+ T cell central memory")
    ,("lymphocyte", "mucosal invariant T cell (MAIT)")
    ,("lymphocyte", "CD4+ T cell effector memory")
    ,("lymphocyte", "CD8+ T cell effector memory")
    ,("CD4-positive, alpha-beta T cell", "CD4+ cytotoxic lymphocyte")
    ,("natural killer cell", "NK CD56bright")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $newSheet.Cells.Item($row, 1).Value = $data[$i][0]
    $newSheet.Cells.Item($row, 2).Value = $data[$i][1]
}

# Apply the header style (bold, centered, bordered) matching the other sheets
$styleSource.Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)

$newSheet.Range("A1").Select()
